$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# CU-01 "CRU promoción": estado pasa de "vacío" a "planificado" y el esfuerzo pasa de 0 a 2
$ws.Range("E5").Value = "planificado"
$ws.Range("F5").Value = 2

# CU-04 "Consultar proximos pagos de alumnos": estado pasa de "vacio" a "planificado" y el esfuerzo pasa de 0 a 2
$ws.Range("E8").Value = "planificado"
$ws.Range("F8").Value = 2

# Actualiza la celda seleccionada/vista de la hoja
$ws.Activate() | Out-Null
$ws.Range("C10").Select() | Out-Null
